$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values: E2:E8 2400 -> 2160, F2:F8 240 -> 336
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 5).Value = 2160
    $ws.Cells.Item($row, 6).Value = 336
}

# Update the active cell selection to F12
$ws.Range("F12").Select()
